# PAL_Deterministic Schedule_V2_H1.xlsx -- "Add files via upload" edit
#
# The commit re-uploads the workbook with an across-the-board -2 shift
# applied to every recorded value in column I (rows 2-108 of Sheet1), and
# moves the active selection from N104 to E79. (The workbook-level
# metadata churn visible in the raw XML diff -- the absPath breadcrumb,
# the xr:revisionPtr save-session GUID, and the bookViews window offset --
# is Excel/session bookkeeping that isn't meaningfully reproducible by
# driving the object model, and isn't data content, so it's left alone.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose ORIGINAL column-I value had been clamped to 5.01. The clamp
# erased the true pre-clamp magnitude, so for these rows "subtract 2 from
# the stored value" does not land on the value recorded after the edit --
# set those few cells to the explicit post-edit values instead.
$explicitIValues = @{
    5  = 2.75
    13 = 2.67
    16 = 2.84
    22 = 2.9000000000000004
    26 = 2.92
    28 = 2.8100000000000005
    43 = 2.96
    51 = 2.84
    55 = 2.84
    61 = 2.67
    62 = 2.9000000000000004
    67 = 2.8900000000000006
    70 = 2.66
}

for ($row = 2; $row -le 108; $row++) {
    $cell = $ws.Cells.Item($row, 9)  # column I
    if ($explicitIValues.ContainsKey($row)) {
        $cell.Value = $explicitIValues[$row]
    } else {
        $cell.Value = $cell.Value2 - 2
    }
}

# Move the live selection to match the saved view (was N104).
$ws.Range("E79").Select() | Out-Null
